# "add baseline to calculation"
#
# The DCASE2023_baseline_task7 row (row 10) on both sheets was missing its
# alg_code (column B) value. This fills it in with "Baseline" on both the
# "audio_quality" and "category_fit" sheets, leaving every other cell,
# style and column untouched.

$wb = $excel.ActiveWorkbook

$wsAudio = $wb.Worksheets.Item("audio_quality")
$wsCategory = $wb.Worksheets.Item("category_fit")

# Fill in the missing alg_code for the baseline submission on the
# "audio_quality" sheet and leave the selection sitting on the edited cell.
$wsAudio.Range("B10").Value = "Baseline"
$wsAudio.Range("B10").Select()

# "category_fit" is the sheet that was active/visible last, so activate it
# and perform the same edit there too.
$wsCategory.Activate()
$wsCategory.Range("B10").Value = "Baseline"
$wsCategory.Range("B15").Select()
